{"js": "// Update the title line (date) in the first paragraph of the document.\nconst firstPara = context.document.body.paragraphs.getFirst();\nconst titleRange = firstPara.getRange();\ntitleRange.insertText(\"2024-02-03 Saturday\", Word.InsertLocation.replace);\n\n// Update the answer grid in the table: replace each populated cell's text\n// with the new value while preserving per-run formatting (font/size/\n// alignment) by replacing the text on the paragraph's Range rather than\n// the cell Body (which would drop paragraph/run formatting).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// (rowIndex, colIndex) -> new text, in document/diff order.\nconst updates = [\n  [0, 0, \"40\u00f79=4, 4\"],\n  [0, 1, \"79\u00f78=9, 7\"],\n  [0, 2, \"45\u00f76=7, 3\"],\n  [0, 3, \"47\u00f76=7, 5\"],\n  [0, 4, \"20\u00f79=2, 2\"],\n  [4, 0, \"81\u00f79=9, 0\"],\n  [4, 1, \"32\u00f77=4, 4\"],\n  [4, 2, \"89\u00f73=29, 2\"],\n  [4, 3, \"81\u00f79=9, 0\"],\n  [4, 4, \"23\u00f74=5, 3\"],\n  [8, 0, \"60\u00f76=10, 0\"],\n  [8, 1, \"93\u00f72=46, 1\"],\n  [8, 2, \"77\u00f79=8, 5\"],\n  [8, 3, \"62\u00f79=6, 8\"],\n  [8, 4, \"55\u00f78=6, 7\"],\n  [12, 0, \"70\u00f74=17, 2\"],\n  [12, 1, \"48\u00f75=9, 3\"],\n  [12, 2, \"80\u00f76=13, 2\"],\n  [12, 3, \"60\u00f72=30, 0\"],\n  [12, 4, \"59\u00f78=7, 3\"],\n  [16, 0, \"48\u00f78=6, 0\"],\n  [16, 1, \"97\u00f75=19, 2\"],\n  [16, 2, \"32\u00f78=4, 0\"],\n  [16, 3, \"30\u00f73=10, 0\"],\n  [16, 4, \"29\u00f74=7, 1\"],\n];\n\nfor (const [rowIndex, colIndex, newText] of updates) {\n  const cell = table.getCell(rowIndex, colIndex);\n  cell.body.paragraphs.load(\"items\");\n  await context.sync();\n\n  const para = cell.body.paragraphs.items[0];\n  const range = para.getRange();\n  range.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the title line (date) in the first paragraph of the document.\n$d.Paragraphs.Item(1).Range.Text = \"2024-02-03 Saturday\"\n\n# Update the answer grid in the table: set each populated cell's Range.Text\n# directly so the existing paragraph/run formatting (font, size, alignment)\n# on that cell is preserved.\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text = \"40\u00f79=4, 4\"\n$t.Cell(1, 2).Range.Text = \"79\u00f78=9, 7\"\n$t.Cell(1, 3).Range.Text = \"45\u00f76=7, 3\"\n$t.Cell(1, 4).Range.Text = \"47\u00f76=7, 5\"\n$t.Cell(1, 5).Range.Text = \"20\u00f79=2, 2\"\n\n$t.Cell(5, 1).Range.Text = \"81\u00f79=9, 0\"\n$t.Cell(5, 2).Range.Text = \"32\u00f77=4, 4\"\n$t.Cell(5, 3).Range.Text = \"89\u00f73=29, 2\"\n$t.Cell(5, 4).Range.Text = \"81\u00f79=9, 0\"\n$t.Cell(5, 5).Range.Text = \"23\u00f74=5, 3\"\n\n$t.Cell(9, 1).Range.Text = \"60\u00f76=10, 0\"\n$t.Cell(9, 2).Range.Text = \"93\u00f72=46, 1\"\n$t.Cell(9, 3).Range.Text = \"77\u00f79=8, 5\"\n$t.Cell(9, 4).Range.Text = \"62\u00f79=6, 8\"\n$t.Cell(9, 5).Range.Text = \"55\u00f78=6, 7\"\n\n$t.Cell(13, 1).Range.Text = \"70\u00f74=17, 2\"\n$t.Cell(13, 2).Range.Text = \"48\u00f75=9, 3\"\n$t.Cell(13, 3).Range.Text = \"80\u00f76=13, 2\"\n$t.Cell(13, 4).Range.Text = \"60\u00f72=30, 0\"\n$t.Cell(13, 5).Range.Text = \"59\u00f78=7, 3\"\n\n$t.Cell(17, 1).Range.Text = \"48\u00f78=6, 0\"\n$t.Cell(17, 2).Range.Text = \"97\u00f75=19, 2\"\n$t.Cell(17, 3).Range.Text = \"32\u00f78=4, 0\"\n$t.Cell(17, 4).Range.Text = \"30\u00f73=10, 0\"\n$t.Cell(17, 5).Range.Text = \"29\u00f74=7, 1\"\n"}
